$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Add Panels")

# Rename the "loading detail" labels used in column I / J (shared strings 27 & 28)
# "Battery Alarm (A)"  -> "Alarm Current(A)"
# "Battery Standby (A)" -> "Standby Current(A)"
$ws.Range("I8:I10").Value = "Alarm Current(A)"
$ws.Range("J8:J10").Value = "Standby Current(A)"

# Update expected default battery standby / alarm-limit values
$ws.Range("F8").Value = 0.39
$ws.Range("F9").Value = 0.227
$ws.Range("F10").Value = 0.251
$ws.Range("G10").Value = 0.439

# Column J now holds wider text ("Standby Current(A)"), so it picks up an
# explicit best-fit width, same as the other label columns.
$ws.Columns.Item(10).ColumnWidth = 18.83

# Restore the cursor/selection state left behind by the edit.
$ws.Range("G9").Select() | Out-Null
